# Updates cryptos list values to match the scraped snapshot.
# Values are prefixed with a literal apostrophe and the cell style is reset
# to "Normal" afterwards so Excel stores them as plain text (matching the
# original inlineStr cells) instead of auto-coercing numeric-looking strings
# (e.g. "1.001", "0.2933") into floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.098.38"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.13%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.898.52"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.27%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'324.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'0.9995"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.36%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4601"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.23%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3882"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.16%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07862"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.88%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9884"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.15%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  -1.52%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.867.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.90%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'5.775"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.22%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'7.037"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -0.57%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("E15").Value = "'  +0.92%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'87.77"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -0.58%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.002"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.27%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000009912"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.49%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'16.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.71%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.9998"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -0.26%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'29.111.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.11%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.319"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.73%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'11.11"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +0.30%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'2.117.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.01%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'2.097"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +2.08%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'156.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -0.25%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'19.43"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.04%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'5.904"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -3.30%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'118.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.40%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  -6.23%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.09324"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.61%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -3.59%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D34").Value = "'1.318"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -2.23%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.133"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -4.24%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'0.05781"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.70%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  -2.51%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = "'0.02083"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  -0.85%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'0.9987"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.38%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'7.662"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -3.40%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.5674"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -1.20%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1794"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -0.22%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'9.684"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -2.30%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.16%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'2.214"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.46%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.5349"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -1.23%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.07004"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.97%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.845"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -1.74%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'2.549"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.18%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'112.72"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.32%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("B51").Value = "'WOONetwork"
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").Value = "'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").Value = "'0.2933"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.29%  "
$ws.Range("E51").Style = "Normal"
